# Edit script: reorder slides and insert one new slide.
#
# Final slide order (by original content):
#   1. "CSV란?" slide              (was slide 3)
#   2. NEW "CSV 파일을 활용해서..." slide
#   3. "리스트 사용" slide          (was slide 5)
#   4. (image-only slide, no text)  (was slide 2)
#   5. "조건문은..." slide          (was slide 4)
#   6. (image-only slide, no text)  (was slide 1)

$p = $ppt.ActivePresentation

function Get-SlideById($id) {
    for ($i = 1; $i -le $p.Slides.Count; $i++) {
        if ($p.Slides.Item($i).SlideID -eq $id) {
            return $p.Slides.Item($i)
        }
    }
    return $null
}

# Remember the SlideIDs of the 5 existing slides before we touch anything.
$origIds = @()
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $origIds += $p.Slides.Item($i).SlideID
}
$idCsvWhat   = $origIds[2]   # was slide 3 - "CSV란?"
$idList      = $origIds[4]   # was slide 5 - "리스트 사용"
$idImgOnlyA  = $origIds[1]   # was slide 2 - image-only
$idCond      = $origIds[3]   # was slide 4 - "조건문은..."
$idImgOnlyB  = $origIds[0]   # was slide 1 - image-only

# 1) Add the new slide at the end, using the same "제목 및 내용" layout as the
#    other text-bearing content slides.
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)
$newId = $newSlide.SlideID

# Remove the inherited (empty) title/content placeholders - the authored
# slide only contains a single free-standing text box.
for ($i = $newSlide.Shapes.Count; $i -ge 1; $i--) {
    $newSlide.Shapes.Item($i).Delete()
}

# 2) Add the text box with the slide's body copy.
$tb = $newSlide.Shapes.AddTextbox(1, 68.07, 55.64, 375.71, 247.19)
$tb.Name = "TextBox 2"

$bodyLines = @(
    "CSV 파일을 활용해서 데이터를 모으고 Python을 활용하여 데이터 정제 및 시각화",
    "",
    "사용 데이터",
    "대구시 연령별 미혼자 수 – 통계청",
    "남자 키 – 병무청",
    "여자 키 – 통계청",
    "남녀 연봉 – 통계청",
    "남녀 운동 여부 – 통계청",
    "남녀 흡연 여부 – 통계청 ",
    "남녀 음주 여부 - 통계청"
)
$tb.TextFrame.TextRange.Text = [string]::Join([char]13, $bodyLines)
$tb.TextFrame.WordWrap = $true
$tb.TextFrame.AutoSize = 1

# 3) Reorder all six slides into their final positions.
$targetOrder = @($idCsvWhat, $newId, $idList, $idImgOnlyA, $idCond, $idImgOnlyB)
for ($pos = 1; $pos -le $targetOrder.Count; $pos++) {
    $slide = Get-SlideById $targetOrder[$pos - 1]
    if ($slide.SlideIndex -ne $pos) {
        $slide.MoveTo($pos)
    }
}
